# "Generate Report for Handback"
#
# This script updates the localization-status workbook to reflect that the
# handback has completed and is now in sync with en-US:
#   - Status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" on all three sheets.
#   - The "Latest Handback DateTime" is refreshed with a new timestamp on the
#     zh-cn and de-de sheets.
#   - The (now resolved) handback-version-mismatch error detail is cleared on
#     the zh-cn and de-de sheets.
#   - A couple of columns are resized to fit the new content.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Overview sheet: Status columns for zh-cn (E2) and de-de (F2) ---
$ws1.Range("E2").Value = $statusText
$ws1.Range("F2").Value = $statusText

# --- zh-cn sheet ---
$ws2.Range("C2").Value = $statusText
$ws2.Range("K2").Value = "2016-08-18 16:47:57"
$ws2.Range("P2").Value = ""

# --- de-de sheet ---
$ws3.Range("C2").Value = $statusText
$ws3.Range("K2").Value = "2016-08-18 16:48:12"
$ws3.Range("P2").Value = ""

# --- Column width adjustments (ColumnWidth is quantized by the engine to the
#     nearest 1/6 character, so we pick the ColumnWidth input that rounds to
#     the desired stored width) ---
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666668   # -> stored width 30
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666668   # -> stored width 30

$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668    # -> stored width 30
$ws2.Columns.Item(16).ColumnWidth = 12.833333333333334   # -> stored width 13.666666666666666

$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668    # -> stored width 30
$ws3.Columns.Item(16).ColumnWidth = 12.833333333333334   # -> stored width 13.666666666666666
